$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.900.93'
$ws.Range('E2').Value = '  +1.83%  '

# Row 3
$ws.Range('D3').Value = '2.621.32'
$ws.Range('E3').Value = '  +1.49%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = "'" + '595.59'
$ws.Range('E5').Value = '  +0.86%  '

# Row 6
$ws.Range('D6').Value = "'" + '155.49'
$ws.Range('E6').Value = '  +0.70%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('E8').Value = '  +1.47%  '

# Row 9
$ws.Range('D9').Value = '2.618.06'
$ws.Range('E9').Value = '  +1.38%  '

# Row 10
$ws.Range('D10').Value = "'" + '0.127'
$ws.Range('E10').Value = '  +10.93%  '

# Row 11
$ws.Range('E11').Value = '  +1.00%  '

# Row 12
$ws.Range('D12').Value = "'" + '5.25'
$ws.Range('E12').Value = '  +0.39%  '

# Row 13
$ws.Range('D13').Value = "'" + '0.355'
$ws.Range('E13').Value = '  -1.37%  '

# Row 14
$ws.Range('D14').Value = "'" + '27.70'
$ws.Range('E14').Value = '  -1.96%  '

# Row 15
$ws.Range('D15').Value = "'" + '0.0000186'
$ws.Range('E15').Value = '  +3.33%  '

# Row 16
$ws.Range('D16').Value = '3.077.79'
$ws.Range('E16').Value = '  +0.36%  '

# Row 17
$ws.Range('D17').Value = '67.719.41'
$ws.Range('E17').Value = '  +1.67%  '

# Row 18
$ws.Range('D18').Value = '2.615.19'
$ws.Range('E18').Value = '  +1.33%  '

# Row 19
$ws.Range('D19').Value = "'" + '11.22'
$ws.Range('E19').Value = '  -0.14%  '

# Row 20
$ws.Range('D20').Value = "'" + '366.59'
$ws.Range('E20').Value = '  +3.09%  '

# Row 21
$ws.Range('D21').Value = "'" + '7.65'
$ws.Range('E21').Value = '  -2.01%  '

# Row 22
$ws.Range('E22').Value = '  -0.70%  '

# Row 23
$ws.Range('D23').Value = "'" + '2.03'
$ws.Range('E23').Value = '  -2.19%  '

# Row 24
$ws.Range('E24').Value = '  -0.07%  '

# Row 25
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').Value = "'" + '9.83'
$ws.Range('E25').Value = '  -6.92%  '

# Row 26
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = "'" + '67.43'
$ws.Range('E26').Value = '  +0.41%  '

# Row 27
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = "'" + '0.0000104'
$ws.Range('E27').Value = '  +0.99%  '

# Row 28
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.744.26'
$ws.Range('E28').Value = '  +0.62%  '

# Row 29
$ws.Range('D29').Value = "'" + '577.52'
$ws.Range('E29').Value = '  -5.91%  '

# Row 30
$ws.Range('D30').Value = "'" + '1.00'
$ws.Range('E30').Value = '  -0.03%  '

# Row 31
$ws.Range('D31').Value = "'" + '1.44'
$ws.Range('E31').Value = '  -1.64%  '

# Row 32
$ws.Range('D32').Value = "'" + '7.93'
$ws.Range('E32').Value = '  -1.01%  '

# Row 33
$ws.Range('E33').Value = '  +0.40%  '

# Row 34
$ws.Range('E34').Value = '  -1.04%  '

# Row 35
$ws.Range('E35').Value = '  -0.08%  '

# Row 36
$ws.Range('D36').Value = "'" + '1.53'
$ws.Range('E36').Value = '  -3.38%  '

# Row 37
$ws.Range('D37').Value = "'" + '4.95'
$ws.Range('E37').Value = '  -2.16%  '

# Row 38
$ws.Range('D38').Value = "'" + '159.07'
$ws.Range('E38').Value = '  +2.75%  '

# Row 39
$ws.Range('D39').Value = "'" + '19.35'
$ws.Range('E39').Value = '  +0.79%  '

# Row 40
$ws.Range('E40').Value = '  -0.05%  '

# Row 41
$ws.Range('D41').Value = "'" + '5.35'
$ws.Range('E41').Value = '  -2.77%  '

# Row 42
$ws.Range('E42').Value = '  +1.35%  '

# Row 43
$ws.Range('E43').Value = '  -3.78%  '

# Row 44
$ws.Range('E44').Value = '  -0.95%  '

# Row 45
$ws.Range('E45').Value = '  +0.05%  '

# Row 46
$ws.Range('E46').Value = '  -0.14%  '

# Row 47
$ws.Range('D47').Value = "'" + '156.11'
$ws.Range('E47').Value = '  +0.10%  '

# Row 48
$ws.Range('E48').Value = '  -6.14%  '

# Row 49
$ws.Range('E49').Value = '  -0.34%  '

# Row 50
$ws.Range('D50').Value = "'" + '0.629'
$ws.Range('E50').Value = '  +3.15%  '

# Row 51
$ws.Range('D51').Value = "'" + '20.86'
$ws.Range('E51').Value = '  -3.10%  '
